# Insert the "datasetSpecializationId" header column before "domain" on the
# SDTMGroup and SDTMGroup1 worksheets, shifting "domain" and "shortName"
# one column to the right (C -> D -> E).

$wb = $excel.ActiveWorkbook

$sheetNames = @("SDTMGroup", "SDTMGroup1")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Shift existing values right: E1 <- D1, D1 <- C1
    $ws.Cells.Item(1, 5).Value2 = $ws.Cells.Item(1, 4).Value2
    $ws.Cells.Item(1, 4).Value2 = $ws.Cells.Item(1, 3).Value2

    # Insert the new header in column C
    $ws.Cells.Item(1, 3).Value2 = "datasetSpecializationId"
}
